# Daily attendance processing - 2025-11-24 04:18:30
#
# Normalizes the ordering of names/emails inside the "Recorded By" (column G)
# cells of the "Session Analysis Results" sheet. The system-recorded entries
# that previously listed the real user after "System"/"system" are reordered
# so the real user (or the more specific entry) comes first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$headerRow = 1
$recordedByCol = 7  # Column G = "Recorded By"

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
